$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.123.00'
$ws.Range('E2').Value = '  +3.51%  '

$ws.Range('D3').Value = '2.246.40'
$ws.Range('E3').Value = '  +0.99%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '231.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.40%  '

$ws.Range('E6').Value = '  -0.47%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '60.96'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.410'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.07%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0913'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.73%  '

$ws.Range('E11').Value = '  +0.53%  '

$ws.Range('D12').Value = '2.579.85'
$ws.Range('E12').Value = '  +1.05%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.62'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.14%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.50'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.66%  '

$ws.Range('E15').Value = '  +2.34%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.804'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.85%  '

$ws.Range('D17').Value = '2.249.05'
$ws.Range('E17').Value = '  +1.15%  '

$ws.Range('D18').Value = '42.970.69'
$ws.Range('E18').Value = '  +3.51%  '

$ws.Range('D19').Value = '0.0₃0938'

$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.21'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.91%  '

$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.25%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '246.32'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.36%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.59'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +8.43%  '

$ws.Range('E24').Value = '  -0.10%  '

$ws.Range('E25').Value = '  +2.89%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.77'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.20%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '169.85'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.07%  '

$ws.Range('E28').Value = '  +1.76%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.49'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.87%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '20.48'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.73%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.69'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.28%  '

$ws.Range('E32').Value = '  -1.13%  '

$ws.Range('E33').Value = '  +1.54%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.97%  '

$ws.Range('E35').Value = '  +5.49%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.44'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.82%  '

$ws.Range('E37').Value = '  +1.64%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.59'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.14%  '

$ws.Range('E39').Value = '  +4.60%  '

$ws.Range('E40').Value = '  +0.15%  '

$ws.Range('E41').Value = '  +0.91%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.000226'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.63%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0973'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.70%  '

$ws.Range('E44').Value = '  -0.52%  '

$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.41'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -9.56%  '

$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '96.93'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.97%  '

$ws.Range('D47').Value = '1.461.94'
$ws.Range('E47').Value = '  -0.13%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '16.43'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.11%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.76'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.72%  '

$ws.Range('E50').Value = '  -0.06%  '

$ws.Range('E51').Value = '  +5.93%  '
